$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 8.931699999999999
$ws.Range("B6").Value = 6.928400000000002
$ws.Range("B7").Value = 5.712000000000002
$ws.Range("D7").Value = -7.582899999999998
$ws.Range("B8").Value = 6.3646
$ws.Range("D11").Value = -7.977599999999998
$ws.Range("D12").Value = -6.4941
$ws.Range("E12").Value = 18.71620000000003
$ws.Range("E13").Value = 16.72420000000001
$ws.Range("E14").Value = 17.0686
$ws.Range("D15").Value = -8.480599999999992
$ws.Range("B16").Value = 5.241399999999999
$ws.Range("E16").Value = 16.4618
$ws.Range("E19").Value = 16.27979999999999
$ws.Range("B20").Value = 9.125899999999987
$ws.Range("D20").Value = -7.863999999999997
$ws.Range("E20").Value = 16.45
$ws.Range("B21").Value = 9.1068
$ws.Range("D21").Value = -8.084500000000002
$ws.Range("D22").Value = -7.7773
$ws.Range("E22").Value = 16.4612
$ws.Range("D23").Value = -7.160299999999998
$ws.Range("B28").Value = 5.890499999999999
$ws.Range("B29").Value = 5.198700000000005
$ws.Range("D29").Value = -7.182399999999994
$ws.Range("B30").Value = 5.065100000000001
$ws.Range("B32").Value = 7.244299999999996
$ws.Range("D34").Value = -7.789600000000002
$ws.Range("E36").Value = 15.8309
$ws.Range("B40").Value = 9.336399999999999
$ws.Range("D42").Value = -8.2066
$ws.Range("D43").Value = -8.081299999999997
$ws.Range("E43").Value = 16.8651
$ws.Range("D44").Value = -8.075199999999995
$ws.Range("D45").Value = -7.726299999999994
$ws.Range("B46").Value = 6.279700000000004
$ws.Range("D46").Value = -8.270499999999997
$ws.Range("E46").Value = 16.6537
$ws.Range("D50").Value = -8.139299999999999
$ws.Range("E50").Value = 16.74499999999999
$ws.Range("B51").Value = 5.445900000000002
$ws.Range("D51").Value = -7.692599999999998
$ws.Range("B52").Value = 5.432600000000001
$ws.Range("B57").Value = 6.022999999999993
$ws.Range("D57").Value = -7.889099999999997
$ws.Range("B59").Value = 5.588799999999996
$ws.Range("B62").Value = 5.568599999999996
$ws.Range("D65").Value = -7.8114
$ws.Range("B66").Value = 5.319099999999996
$ws.Range("D66").Value = -7.317900000000002
$ws.Range("D67").Value = -6.428300000000003
$ws.Range("B73").Value = 8.699499999999999
$ws.Range("B74").Value = 9.204899999999995
$ws.Range("E76").Value = 16.2411
$ws.Range("B77").Value = 8.65140000000001
$ws.Range("D79").Value = -6.224500000000003
$ws.Range("D84").Value = -9.313099999999999
$ws.Range("D87").Value = -8.088099999999995
$ws.Range("B92").Value = 4.724800000000003
$ws.Range("D92").Value = -6.467
$ws.Range("E95").Value = 17.97680000000001
$ws.Range("D97").Value = -7.846199999999993
$ws.Range("E97").Value = 16.93660000000001
$ws.Range("E99").Value = 16.34659999999999
$ws.Range("B100").Value = 5.777499999999996
